$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 5 to make room for the new "ECs" sending-cluster block
# (shifts the old MuSCs block from rows 5-7 down to rows 8-10)
$ws.Rows("5:7").Insert()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cntn2"
$ws.Range("C2").Value = "Nrcam"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.05028033333333334
$ws.Range("H2").Value = 0.150841
$ws.Range("I2").Value = 0.1252715694221136
$ws.Range("J2").Value = 0.1252715694221136
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.027105
$ws.Range("N2").Value = 0.081315
$ws.Range("O2").Value = 0.01418818755522427
$ws.Range("P2").Value = 0.01418818755522427
$ws.Range("Q2").Value = 0.001362848435
$ws.Range("R2").Value = 0.012265635915
$ws.Range("S2").Value = 0.001777376522298245
$ws.Range("T2").Value = 0.001777376522298245

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cntn2"
$ws.Range("C3").Value = "Nrcam"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.05028033333333334
$ws.Range("H3").Value = 0.150841
$ws.Range("I3").Value = 0.1252715694221136
$ws.Range("J3").Value = 0.1252715694221136
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.4031683333333334
$ws.Range("N3").Value = 1.209505
$ws.Range("O3").Value = 0.2110395841970304
$ws.Range("P3").Value = 0.2110395841970304
$ws.Range("Q3").Value = 0.02027143818944445
$ws.Range("R3").Value = 0.182442943705
$ws.Range("S3").Value = 0.02643725992255229
$ws.Range("T3").Value = 0.02643725992255228

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cntn2"
$ws.Range("C4").Value = "Nrcam"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.05028033333333334
$ws.Range("H4").Value = 0.150841
$ws.Range("I4").Value = 0.1252715694221136
$ws.Range("J4").Value = 0.1252715694221136
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.480118666666667
$ws.Range("N4").Value = 4.440356
$ws.Range("O4").Value = 0.7747722282477453
$ws.Range("P4").Value = 0.7747722282477453
$ws.Range("Q4").Value = 0.07442085993288888
$ws.Range("R4").Value = 0.6697877393959999
$ws.Range("S4").Value = 0.09705693297726306
$ws.Range("T4").Value = 0.09705693297726303

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cntn2"
$ws.Range("C5").Value = "Nrcam"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.300673
$ws.Range("H5").Value = 0.9020189999999999
$ws.Range("I5").Value = 0.7491155307811899
$ws.Range("J5").Value = 0.7491155307811898
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.027105
$ws.Range("N5").Value = 0.081315
$ws.Range("O5").Value = 0.01418818755522427
$ws.Range("P5").Value = 0.01418818755522427
$ws.Range("Q5").Value = 0.008149741665
$ws.Range("R5").Value = 0.07334767498499999
$ws.Range("S5").Value = 0.0106285916512549
$ws.Range("T5").Value = 0.0106285916512549

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cntn2"
$ws.Range("C6").Value = "Nrcam"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.300673
$ws.Range("H6").Value = 0.9020189999999999
$ws.Range("I6").Value = 0.7491155307811899
$ws.Range("J6").Value = 0.7491155307811898
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.4031683333333334
$ws.Range("N6").Value = 1.209505
$ws.Range("O6").Value = 0.2110395841970304
$ws.Range("P6").Value = 0.2110395841970304
$ws.Range("Q6").Value = 0.1212218322883333
$ws.Range("R6").Value = 1.090996490595
$ws.Range("S6").Value = 0.1580930301316001
$ws.Range("T6").Value = 0.1580930301316001

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cntn2"
$ws.Range("C7").Value = "Nrcam"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.300673
$ws.Range("H7").Value = 0.9020189999999999
$ws.Range("I7").Value = 0.7491155307811899
$ws.Range("J7").Value = 0.7491155307811898
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.480118666666667
$ws.Range("N7").Value = 4.440356
$ws.Range("O7").Value = 0.7747722282477453
$ws.Range("P7").Value = 0.7747722282477453
$ws.Range("Q7").Value = 0.4450317198626666
$ws.Range("R7").Value = 4.005285478763999
$ws.Range("S7").Value = 0.5803939089983349
$ws.Range("T7").Value = 0.5803939089983348

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Cntn2"
$ws.Range("C8").Value = "Nrcam"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.05041733333333333
$ws.Range("H8").Value = 0.151252
$ws.Range("I8").Value = 0.1256128997966967
$ws.Range("J8").Value = 0.1256128997966967
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.027105
$ws.Range("N8").Value = 0.081315
$ws.Range("O8").Value = 0.01418818755522427
$ws.Range("P8").Value = 0.01418818755522427
$ws.Range("Q8").Value = 0.00136656182
$ws.Range("R8").Value = 0.01229905638
$ws.Range("S8").Value = 0.001782219381671125
$ws.Range("T8").Value = 0.001782219381671125

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Cntn2"
$ws.Range("C9").Value = "Nrcam"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.05041733333333333
$ws.Range("H9").Value = 0.151252
$ws.Range("I9").Value = 0.1256128997966967
$ws.Range("J9").Value = 0.1256128997966967
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4031683333333334
$ws.Range("N9").Value = 1.209505
$ws.Range("O9").Value = 0.2110395841970304
$ws.Range("P9").Value = 0.2110395841970304
$ws.Range("Q9").Value = 0.02032667225111111
$ws.Range("R9").Value = 0.18294005026
$ws.Range("S9").Value = 0.02650929414287812
$ws.Range("T9").Value = 0.02650929414287811

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Cntn2"
$ws.Range("C10").Value = "Nrcam"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.05041733333333333
$ws.Range("H10").Value = 0.151252
$ws.Range("I10").Value = 0.1256128997966967
$ws.Range("J10").Value = 0.1256128997966967
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.480118666666667
$ws.Range("N10").Value = 4.440356
$ws.Range("O10").Value = 0.7747722282477453
$ws.Range("P10").Value = 0.7747722282477453
$ws.Range("Q10").Value = 0.07462363619022222
$ws.Range("R10").Value = 0.6716127257119999
$ws.Range("S10").Value = 0.09732138627214744
$ws.Range("T10").Value = 0.09732138627214741
